$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reverse-Block($ws, $startRow, $endRow) {
    $rng = $ws.Range("E$startRow`:F$endRow")
    $vals = $rng.Value2
    $n = $vals.GetLength(0)
    $newvals = New-Object 'object[,]' $n,2
    for ($i=0; $i -lt $n; $i++) {
        $newvals[$i,0] = $vals[$n-$i,1]
        $newvals[$i,1] = $vals[$n-$i,2]
    }
    $rng.Value2 = $newvals
}

Reverse-Block $ws 16 43
Reverse-Block $ws 45 103
